$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "(system)API Search Form Visibility"
$ws.Range("E2:E7").Value = "Y"

$ws.Range("I5").Select() | Out-Null
